# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) previously held a derived "Strike#"-style value.
# It has been recalculated/regenerated using the true strikeout counts (K)
# for each game and is re-written here, row by row, for rows 2-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (strikeout) value for column G
$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 2
    6  = 1
    7  = 2
    8  = 2
    9  = 1
    10 = 1
    11 = 0
    12 = 2
    13 = 0
    14 = 2
    15 = 1
    16 = 0
    17 = 0
    18 = 2
    19 = 1
    20 = 0
    21 = 0
    22 = 2
    23 = 1
    24 = 3
    25 = 2
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 3
    31 = 0
    32 = 0
    33 = 2
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 3
    39 = 1
    40 = 0
    41 = 1
    42 = 3
    43 = 3
    44 = 1
    45 = 1
    46 = 2
    47 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
